$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $value
    $c.Style = "Normal"
}

Set-TextValue 'D2' '27.571.55'
Set-TextValue 'E2' '  -2.67%  '
Set-TextValue 'D3' '1.752.94'
Set-TextValue 'E3' '  -3.56%  '
Set-TextValue 'D4' '1.002'
Set-TextValue 'E4' '  +0.17%  '
Set-TextValue 'D5' '324.38'
Set-TextValue 'E5' '  -0.33%  '
Set-TextValue 'D6' '1.001'
Set-TextValue 'E6' '  +0.15%  '
Set-TextValue 'D7' '0.4482'
Set-TextValue 'E7' '  +3.10%  '
Set-TextValue 'D8' '0.3622'
Set-TextValue 'E8' '  -1.52%  '
Set-TextValue 'E9' '  -2.16%  '
Set-TextValue 'D10' '42.16'
Set-TextValue 'E10' '  -6.07%  '
Set-TextValue 'D11' '1.103'
Set-TextValue 'E11' '  -3.87%  '
Set-TextValue 'D12' '1.001'
Set-TextValue 'E12' '  +0.09%  '
Set-TextValue 'D13' '20.69'
Set-TextValue 'E13' '  -6.09%  '
Set-TextValue 'D14' '6.038'
Set-TextValue 'E14' '  -4.57%  '
Set-TextValue 'D15' '7.175'
Set-TextValue 'E15' '  -4.27%  '
Set-TextValue 'D16' '1.751.69'
Set-TextValue 'E16' '  -3.64%  '
Set-TextValue 'D17' '92.89'
Set-TextValue 'E17' '  -2.59%  '
Set-TextValue 'E18' '  -1.42%  '
Set-TextValue 'D19' '0.06396'
Set-TextValue 'E19' '  -0.75%  '
Set-TextValue 'D20' '1.000'
Set-TextValue 'E20' '  +0.06%  '
Set-TextValue 'D21' '16.91'
Set-TextValue 'E21' '  -3.11%  '
Set-TextValue 'D22' '5.858'
Set-TextValue 'E22' '  -6.20%  '
Set-TextValue 'D23' '27.611.35'
Set-TextValue 'E23' '  -2.56%  '
Set-TextValue 'E24' '  -3.15%  '
Set-TextValue 'D25' '2.109'
Set-TextValue 'E25' '  -1.55%  '
Set-TextValue 'D26' '161.80'
Set-TextValue 'E26' '  +1.23%  '
Set-TextValue 'E27' '  -1.69%  '
Set-TextValue 'D28' '1.952.15'
Set-TextValue 'E28' '  -3.68%  '
Set-TextValue 'E29' '  -6.89%  '
Set-TextValue 'D30' '125.25'
Set-TextValue 'E30' '  -4.92%  '
Set-TextValue 'E31' '  -9.98%  '
Set-TextValue 'D32' '3.664'
Set-TextValue 'E32' '  +2.95%  '
Set-TextValue 'D33' '0.09018'
Set-TextValue 'E33' '  -1.29%  '
Set-TextValue 'D34' '5.551'
Set-TextValue 'E34' '  -7.80%  '
Set-TextValue 'E35' '  -7.75%  '
Set-TextValue 'D36' '0.02310'
Set-TextValue 'E36' '  -4.03%  '
Set-TextValue 'E37' '  -3.71%  '
Set-TextValue 'D38' '0.6384'
Set-TextValue 'E38' '  -3.33%  '
Set-TextValue 'D39' '4.981'
Set-TextValue 'E39' '  -4.68%  '
Set-TextValue 'E40' '  -3.59%  '
Set-TextValue 'D41' '1.193'
Set-TextValue 'E41' '  -0.80%  '
Set-TextValue 'D42' '1.000'
Set-TextValue 'E42' '  +0.13%  '
Set-TextValue 'D43' '1.391'
Set-TextValue 'E43' '  -2.69%  '
Set-TextValue 'D44' '7.808'
Set-TextValue 'E44' '  -2.94%  '
Set-TextValue 'D45' '13.30'
Set-TextValue 'E45' '  -4.24%  '
Set-TextValue 'D46' '0.5906'
Set-TextValue 'E46' '  -3.34%  '
Set-TextValue 'D47' '3.707'
Set-TextValue 'E47' '  -0.78%  '
Set-TextValue 'B48' 'Quant'
Set-TextValue 'C48' 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
Set-TextValue 'D48' '121.90'
Set-TextValue 'E48' '  -3.24%  '
Set-TextValue 'B49' 'NEARProtocol'
Set-TextValue 'C49' 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
Set-TextValue 'D49' '1.958'
Set-TextValue 'E49' '  -3.15%  '
Set-TextValue 'E50' '  -0.47%  '
Set-TextValue 'D51' '0.06874'
Set-TextValue 'E51' '  -1.77%  '
